# Add kanban board for the completed fourth sprint
# Updates the ARCHIVE sheet of the Backlog workbook: replaces the "assigned
# to" placeholder, re-shuffles a handful of rows at the bottom of the board
# and introduces 4 brand-new backlog items for sprint 4.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ARCHIVE")
$ws.Activate()

# --- Full data block (rows 3-48, columns B:G) -----------------------------
# Column order: B=CATEGORY, C=ASSIGNED TO, D=ACTION, E=STATUS, F=PRIORITY, G=HOURS
$data = @(
    @("Task",     "Ettore Farinelli", "Initialize project repository", "Complete", "High", 0.5),
    @("Task",     "Ettore Farinelli", "Set up build system", "Complete", "High", 0.5),
    @("Research", "Ettore Farinelli", "Configure continuous integration (CI) and automated testing", "Complete", "Low", 0.5),
    @("Task",     "Ettore Farinelli", "Set up code formatting, linting, and .gitignore", "Complete", "Medium", 0.5),
    @("Content",  "Ettore Farinelli", "Set up documentation system (Jekyll + Just the Docs)", "Complete", "Low", 0.5),
    @("Research", "Ettore Farinelli", "Define clean architecture folder structure (domain, application, infrastructure)", "Complete", "High", 0.5),
    @("Task",     "Ettore Farinelli", "Define Backlog for the entire project", "Complete", "High", 1),
    @("Feature",  "Ettore Farinelli", "Create initial project documentation and architecture overview", "Complete", "High", 1),
    @("Feature",  "Ettore Farinelli", "Design and document core game loop", "Complete", "High", 0.5),
    @("Feature",  "Ettore Farinelli", "Define main entry point (e.g., CliApp.scala)", "Complete", "High", 0.5),
    @("Feature",  "Ettore Farinelli", "Implement error handling and logging infrastructure", "Complete", "Medium", 1),
    @("Feature",  "Ettore Farinelli", "Set up configuration management (game settings, environment variables)", "Complete", "Low", 0.5),
    @("Feature",  "Ettore Farinelli", "Define Player model and attributes", "Complete", "High", 1.5),
    @("Feature",  "Ettore Farinelli", "Define World model and structure", "Complete", "High", 1.5),
    @("Feature",  "Ettore Farinelli", "Define GameState management", "Complete", "High", 1),
    @("Feature",  "Ettore Farinelli", "Specify interfaces/traits for world generation", "Complete", "Medium", 1),
    @("Feature",  "Ettore Farinelli", "Specify interfaces/traits for player actions and game events", "Complete", "Medium", 1.5),
    @("Feature",  "Ettore Farinelli", "Implement core domain logic for turn progression, win/loss conditions", "Complete", "High", 1.5),
    @("Feature",  "Ettore Farinelli", "Implement CLI controller for user input/output", "Complete", "Medium", 1.5),
    @("Feature",  "Ettore Farinelli", "Implement main menu navigation (start, help, quit)", "Complete", "Low", 1.5),
    @("Feature",  "Ettore Farinelli", "Implement command parsing and validation", "Complete", "Medium", 1),
    @("Update",   "Ettore Farinelli", "Connect CLI controller to application layer (use case invocation)", "Complete", "Medium", 1.5),
    @("Research", "Ettore Farinelli", "Add input/output adapters for future UI expansion (e.g., GUI, web)", "Complete", "Medium", 1),
    @("Feature",  "Ettore Farinelli", "Disaply world", "Complete", "Low", 2),
    @("Feature",  "Ettore Farinelli", "Implement use cases for player actions (move, interact, inventory, etc.)", "Complete", "Medium", 1.5),
    @("Feature",  "Ettore Farinelli", "Implement use cases for world generation and initialization", "Complete", "Medium", 2),
    @("Feature",  "Ettore Farinelli", "Design procedural world generation algorithms (noise, cellular automata, etc.)", "Complete", "Medium", 3),
    @("Feature",  "Ettore Farinelli", "Implement terrain generation (map, tiles, biomes)", "Complete", "Medium", 1),
    @("Feature",  "Ettore Farinelli", "Implement structure/object placement (buildings, dungeons, items)", "Complete", "Medium", 3),
    @("Task",     "Ettore Farinelli", "Write tests for world generation consistency", "Complete", "Medium", 1),
    @("Feature",  "Ettore Farinelli", "Implement player movement and interaction", "Complete", "Medium", 2),
    @("Feature",  "Ettore Farinelli", "Implement inventory system and item management", "Complete", "Medium", 3),
    @("Feature",  "Ettore Farinelli", "Implement Enemy  movement system in prolog", "Complete", "Medium", 6),
    @("Feature",  "Ettore Farinelli", "Implement Enemy  combact system", "Complete", "Medium", 3),
    @("Feature",  "Ettore Farinelli", "Implement combat or challenge mechanics ", "Complete", "Medium", 2),
    @("Feature",  "Ettore Farinelli", "Refactor view to use monads", "Complete", "Low", 3),
    @("Feature",  "Ettore Farinelli", "Implement state history with Monad", "Complete", "Low", 3),
    @("Feature",  "Ettore Farinelli", "Implement scoring, progression, or achievement systems", "Complete", "Low", 1),
    @("Research", "Ettore Farinelli", "Set up and maintain automated unit tests for all major systems", "In Progress", "Medium", 3),
    @("Research", "Ettore Farinelli", "Set up integration tests for game loop and world generation", "In Progress", "Low", 2),
    @("Task",     "Ettore Farinelli", "Perform manual playtesting and bug tracking", "Complete", "Low", 1.5),
    @("Task",     "Ettore Farinelli", "Maintain up-to-date developer documentation (architecture, API, setup)", "In Progress", "Low", $null),
    @("Update",   "Ettore Farinelli", "Update user/player documentation (how to play, features)", "In Progress", "Low", $null),
    @("Content",  "Ettore Farinelli", "Add diagrams and visual aids to docs (architecture, workflows)", "In Progress", "Low", $null),
    @("Update",   "Ettore Farinelli", "Refine CLI or UI for usability and clarity", "Complete", "Low", 0.5),
    @("Update",   "Ettore Farinelli", "Add error messages and help screens", "Complete", "Low", 1)
)

# Shared-string allocation in the engine follows *write order* for brand-new
# strings. The canonical OOXML introduces "...combact system" (row 36) ahead
# of "...movement system in prolog" (row 35) in the shared string table, even
# though row 35 is rendered first, so pre-seed D36 here to claim that slot.
$ws.Cells.Item(36, 4).Value = "Implement Enemy  combact system"

$startRow = 3
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $entry = $data[$i]
    $ws.Cells.Item($row, 2).Value = $entry[0]
    $ws.Cells.Item($row, 3).Value = $entry[1]
    $ws.Cells.Item($row, 4).Value = $entry[2]
    $ws.Cells.Item($row, 5).Value = $entry[3]
    $ws.Cells.Item($row, 6).Value = $entry[4]
    if ($null -eq $entry[5]) {
        $ws.Cells.Item($row, 7).Value = $null
    } else {
        $ws.Cells.Item($row, 7).Value = $entry[5]
    }
}

# --- Row height tweaks (the "odd" 21.6pt row slides from 38 down to 40) ---
$ws.Rows.Item(38).RowHeight = 22.05
$ws.Rows.Item(40).RowHeight = 21.6

# --- View state: scrolled down a couple of rows, selection moved to E46 ---
$ws.Application.ActiveWindow.ScrollRow = 33
$ws.Range("E46").Select()
